# The workbook's single sheet ("quadratic-svm-score") holds a small
# prediction/score table:
#   A: Row label (text)   B: score value (number)   C: prediction (number)
#
# This refresh re-ran the upstream pipeline against an updated copy of
# ful-path.csv, producing new (real) score values for column B instead
# of the old placeholder "1" values. Column A/C and the header row are
# unchanged in content; only the numeric scores move.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1001.6842732210839
$ws.Range("B3").Value = 1501.6868740106941
$ws.Range("B4").Value = 1788.2785062091918
